$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header "NOM" -> "nom"
$ws.Range("A1").Value = "nom"

# Append new rows at the bottom of the table (previously A1:B165, now A1:B168)
$ws.Range("A166").Value = "digigrowing"
$ws.Range("B166").Value = "ceo.digigrowing@gmail.com"

$ws.Range("A167").Value = "tradrly"
$ws.Range("B167").Value = "tradrlyjob@gmail.com"

$ws.Range("A168").Value = "SHOP MY INFLUENCE"
$ws.Range("B168").Value = "contact@influens.fr"

# Match formatting (row height + alignment) of the existing data rows
$ws.Rows.Item(166).RowHeight = 18.75
$ws.Rows.Item(167).RowHeight = 18.75
$ws.Rows.Item(168).RowHeight = 18.75

$ws.Range("A166:B168").HorizontalAlignment = $ws.Range("A165:B165").HorizontalAlignment
